$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 entirely (the invalid "61.076.055/0001" entry).
# This shifts row 10 ("59.981.829/0001-65") up to become the new row 9.
$ws.Rows.Item(9).Delete()

# Reset the active cell/selection back to A1 (the sheet's default selection),
# since the deleted row's selection (A10) is no longer valid.
$ws.Range("A1").Select() | Out-Null
